# edit.ps1 -- apply "New crime data collected" weekly update
# Target workbook: CompStat weekly crime report (9th Precinct)
#
# Summary of the edit:
#   - Report header: issue Number 22 -> 23; week 5/26/2025-6/1/2025 -> 6/2/2025-6/8/2025
#   - Crime-stat grid (rows 14-28, cols C:N) refreshed with this week's figures.
#     Some cells switch between a numeric value and the literal placeholder
#     text "0" / "***.*" (used when a % change is undefined, e.g. division by
#     zero) depending on the new numbers, so both the stored value AND the
#     cell's text/number type need updating to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Donor cells: stable, never-edited cells already carrying the exact style
# we need (so re-using their format via PasteSpecial keeps the same style
# index instead of Excel minting a brand-new one for a one-off NumberFormat
# tweak).
#   $donorText   -> style used by text cells in the grid ("0" / "***.*")
#   $donorNum14  -> style used by integer-formatted numeric cells
#   $donorNum15  -> style used by decimal-formatted numeric cells (% chg)
# ---------------------------------------------------------------------------
$donorText = $ws.Range("A14")
$donorNum14 = $ws.Range("I29")
$donorNum15 = $ws.Range("K29")

# ---------------------------------------------------------------------------
# Header / report metadata
# ---------------------------------------------------------------------------
$ws.Range("C8").Characters(21, 2).Text = "23"
$ws.Range("C9").Characters(27, 9).Text = "6/2/2025"
$ws.Range("C9").Characters(46, 8).Text = "6/8/2025"

# ---------------------------------------------------------------------------
# Crime statistics grid
# ---------------------------------------------------------------------------
$ws.Range("N14").Value = -80
$ws.Range("C15").NumberFormat = "@"
$ws.Range("C15").Value = "0"
$donorText.Copy()
$ws.Range("C15").PasteSpecial(-4122)
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0"
$donorText.Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "***.*"
$donorText.Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("N15").Value = -66.666666666666
$ws.Range("C16").Value = 3
$ws.Range("D16").Value = 1
$ws.Range("E16").Value = 200
$ws.Range("F16").Value = 14
$ws.Range("G16").Value = 11
$ws.Range("H16").Value = 27.272727272727
$ws.Range("I16").Value = 53
$ws.Range("J16").Value = 55
$ws.Range("K16").Value = -3.636363636363
$ws.Range("L16").Value = -32.051282051282
$ws.Range("M16").Value = -22.058823529411
$ws.Range("N16").Value = -85.195530726257
$donorNum14.Copy()
$ws.Range("C17").PasteSpecial(-4122)
$ws.Range("C17").Value = 3
$ws.Range("D17").Value = 3
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 8
$ws.Range("G17").Value = 13
$ws.Range("H17").Value = -38.461538461538
$ws.Range("I17").Value = 79
$ws.Range("J17").Value = 79
$ws.Range("L17").Value = -26.168224299065
$ws.Range("M17").Value = 17.910447761194
$ws.Range("N17").Value = -65.800865800865
$ws.Range("C18").Value = 2
$ws.Range("D18").Value = 3
$ws.Range("F18").Value = 9
$ws.Range("G18").Value = 18
$ws.Range("H18").Value = -50
$ws.Range("I18").Value = 73
$ws.Range("J18").Value = 75
$ws.Range("K18").Value = -2.666666666666
$ws.Range("L18").Value = -41.129032258064
$ws.Range("M18").Value = -29.807692307692
$ws.Range("N18").Value = -78.338278931750
$ws.Range("C19").Value = 15
$ws.Range("D19").Value = 15
$ws.Range("E19").Value = 0
$ws.Range("F19").Value = 72
$ws.Range("G19").Value = 54
$ws.Range("H19").Value = 33.333333333333
$ws.Range("I19").Value = 358
$ws.Range("J19").Value = 320
$ws.Range("K19").Value = 11.875
$ws.Range("L19").Value = -17.321016166281
$ws.Range("M19").Value = 10.153846153846
$ws.Range("N19").Value = -43.799058084772
$donorNum14.Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("D20").Value = 3
$donorNum15.Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("E20").Value = -100
$ws.Range("F20").Value = 1
$donorNum14.Copy()
$ws.Range("G20").PasteSpecial(-4122)
$ws.Range("G20").Value = 3
$donorNum15.Copy()
$ws.Range("H20").PasteSpecial(-4122)
$ws.Range("H20").Value = -66.666666666666
$ws.Range("J20").Value = 15
$ws.Range("K20").Value = 6.666666666666
$ws.Range("L20").Value = 6.666666666666
$ws.Range("M20").Value = -20
$ws.Range("N20").Value = -92.792792792792
$ws.Range("C21").Value = 23
$ws.Range("E21").Value = -8
$ws.Range("F21").Value = 108
$ws.Range("G21").Value = 100
$ws.Range("H21").Value = 8
$ws.Range("I21").Value = 587
$ws.Range("J21").Value = 553
$ws.Range("K21").Value = 6.148282097649
$ws.Range("L21").Value = -23.368146214099
$ws.Range("M21").Value = -0.508474576271
$ws.Range("N21").Value = -67.586968525676
$ws.Range("C22").NumberFormat = "@"
$ws.Range("C22").Value = "0"
$donorText.Copy()
$ws.Range("C22").PasteSpecial(-4122)
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 10
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 57
$ws.Range("J23").Value = 34
$ws.Range("K23").Value = 67.647058823529
$ws.Range("L23").Value = -8.064516129032
$ws.Range("M23").Value = 18.75
$ws.Range("C24").Value = 24
$ws.Range("D24").Value = 27
$ws.Range("E24").Value = -11.111111111111
$ws.Range("F24").Value = 90
$ws.Range("G24").Value = 122
$ws.Range("H24").Value = -26.229508196721
$ws.Range("I24").Value = 659
$ws.Range("J24").Value = 632
$ws.Range("K24").Value = 4.272151898734
$ws.Range("L24").Value = 6.634304207119
$ws.Range("M24").Value = -4.906204906204
$ws.Range("C25").Value = 7
$ws.Range("D25").Value = 18
$ws.Range("E25").Value = -61.111111111111
$ws.Range("F25").Value = 42
$ws.Range("G25").Value = 74
$ws.Range("H25").Value = -43.243243243243
$ws.Range("I25").Value = 320
$ws.Range("J25").Value = 368
$ws.Range("K25").Value = -13.043478260869
$ws.Range("L25").Value = -12.328767123287
$ws.Range("C26").Value = 10
$ws.Range("D26").Value = 13
$ws.Range("E26").Value = -23.076923076923
$ws.Range("F26").Value = 43
$ws.Range("G26").Value = 33
$ws.Range("H26").Value = 30.303030303030
$ws.Range("I26").Value = 184
$ws.Range("J26").Value = 170
$ws.Range("K26").Value = 8.235294117647
$ws.Range("L26").Value = -11.538461538461
$ws.Range("M26").Value = -5.641025641025
$ws.Range("C27").NumberFormat = "@"
$ws.Range("C27").Value = "0"
$donorText.Copy()
$ws.Range("C27").PasteSpecial(-4122)
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0"
$donorText.Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "***.*"
$donorText.Copy()
$ws.Range("E27").PasteSpecial(-4122)
$donorNum14.Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 6
$ws.Range("G28").Value = 5
$ws.Range("H28").Value = 20
$ws.Range("I28").Value = 23
$ws.Range("K28").Value = -20.689655172413
$ws.Range("L28").Value = 4.545454545454

